$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Write cells in the same first-use order as the source shared strings table
# (filename, EN x2, RU x2, converted x2) so sharedStrings.xml indices line up.
$ws.Range("A2").Value = 'SCRIPT/T01P02A/enter09.ssb'
$ws.Range("C2").Value = ' We all have a purpose. One that\nis unique to each of us.'
$ws.Range("C3").Value = ' One must serve that purpose the\nbest they can. Kwah... Hah!'
$ws.Range("D2").Value = ' У нас всех есть цель в жизни.\nУ каждого она своя.'
$ws.Range("D3").Value = ' Каждый должен достичь этой\nцели, приложив максимум усилий. Квох... Ха!'
$ws.Range("E2").Value = ' Ô îàò âòåö åòóû øåìû â çéèîé.\nÔ ëàçäïãï ïîà òâïÿ.'
$ws.Range("E3").Value = ' Ëàçäúê äïìçåî äïòóéœû üóïê\nøåìé, ðñéìïçéâ íàëòéíôí ôòéìéê. Ëâïö... Öà!'

$ws.Range("B2").Value = 311
$ws.Range("B3").Value = 314

$ws.Rows("2:2").RowHeight = 43.2
$ws.Rows("3:3").RowHeight = 21.6

[void]$ws.Range("D5").Select()
